# Scheduled-runner refresh of market-board price snapshots (currentAveragePrice /
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ -- columns H:N) for a
# handful of leves across the crafting-job sheets. Text columns (A:G) are untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 53: No Accounting for Waste / Enchanted Electrum Ink
$ws.Range("H53").Value = 1740.7222
$ws.Range("I53").Value = 2043.2222
$ws.Range("J53").Value = 1438.2222
$ws.Range("K53").Value = 2043.2222
$ws.Range("L53").Value = 1438.2222
$ws.Range("M53").Value = -1406.2222
$ws.Range("N53").Value = -2712.2222
# Row 69: Steeling the Knife, Steeling the Mind / Grade 1 Mind Dissolvent
$ws.Range("H69").Value = 7577.5557
$ws.Range("J69").Value = 10549.667
$ws.Range("L69").Value = 31649.001
$ws.Range("N69").Value = -33397.001
# Row 72: Surgical Substitution (L) / Grade 1 Mind Dissolvent
$ws.Range("H72").Value = 7577.5557
$ws.Range("J72").Value = 10549.667
$ws.Range("L72").Value = 94947.003
$ws.Range("N72").Value = -103683.003
# Row 80: Cleansing the Wicked Humours / Hallowed Water
$ws.Range("H80").Value = 4863.3335
$ws.Range("I80").Value = 3775.2
$ws.Range("J80").Value = 5640.5713
$ws.Range("K80").Value = 11325.6
$ws.Range("L80").Value = 16921.7139
$ws.Range("M80").Value = -10327.6
$ws.Range("N80").Value = -18917.7139
# Row 83: Washing Away the Sins (L) / Hallowed Water
$ws.Range("H83").Value = 4863.3335
$ws.Range("I83").Value = 3775.2
$ws.Range("J83").Value = 5640.5713
$ws.Range("K83").Value = 33976.8
$ws.Range("L83").Value = 50765.14169999999
$ws.Range("M83").Value = -28984.8
$ws.Range("N83").Value = -60749.14169999999
# Row 98: The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 2779.6
$ws.Range("I98").Value = 2779.6
$ws.Range("K98").Value = 2779.6
$ws.Range("M98").Value = -1281.6
# Row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 2779.6
$ws.Range("I122").Value = 2779.6
$ws.Range("K122").Value = 8338.8
$ws.Range("M122").Value = -5888.799999999999
# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 2099.45
$ws.Range("I132").Value = 1474.2285
$ws.Range("J132").Value = 6476
$ws.Range("K132").Value = 4422.6855
$ws.Range("L132").Value = 19428
$ws.Range("M132").Value = -1892.6855
$ws.Range("N132").Value = -24488
# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 2559.7576
$ws.Range("I137").Value = 1764.4348
$ws.Range("K137").Value = 5293.3044
$ws.Range("M137").Value = -2743.3044
# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 6063278
$ws.Range("I138").Value = 1847.3158
$ws.Range("J138").Value = 14289506
$ws.Range("K138").Value = 5541.9474
$ws.Range("L138").Value = 42868518
$ws.Range("M138").Value = -401.9474
$ws.Range("N138").Value = -42878798

$ws = $wb.Worksheets.Item("ARM")
# Row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Range("H102").Value = 2185.7144
$ws.Range("I102").Value = 2050
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 2050
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -428
$ws.Range("N102").Value = -6244
# Row 103: Sweeping the Legs / Doman Steel Greaves of Striking
$ws.Range("H103").Value = 28120.666
$ws.Range("J103").Value = 28120.666
$ws.Range("L103").Value = 28120.666
$ws.Range("N103").Value = -30464.666
# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 2260.4348
$ws.Range("I122").Value = 1594.2222
$ws.Range("J122").Value = 2688.7144
$ws.Range("K122").Value = 4782.6666
$ws.Range("L122").Value = 8066.1432
$ws.Range("M122").Value = -2332.6666
$ws.Range("N122").Value = -12966.1432

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 1610.15
$ws.Range("I86").Value = 1494.2941
$ws.Range("J86").Value = 2266.6667
$ws.Range("K86").Value = 1494.2941
$ws.Range("L86").Value = 2266.6667
$ws.Range("M86").Value = -371.2941000000001
$ws.Range("N86").Value = -4512.6667
# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 1610.15
$ws.Range("I89").Value = 1494.2941
$ws.Range("J89").Value = 2266.6667
$ws.Range("K89").Value = 7471.4705
$ws.Range("L89").Value = 11333.3335
$ws.Range("M89").Value = -1855.4705
$ws.Range("N89").Value = -22565.3335

$ws = $wb.Worksheets.Item("CRP")
# Row 62: Splinter in the Sewers / Cedar Lumber
$ws.Range("H62").Value = 2728.6843
$ws.Range("I62").Value = 2238.3845
$ws.Range("J62").Value = 3791
$ws.Range("K62").Value = 2238.3845
$ws.Range("L62").Value = 3791
$ws.Range("M62").Value = -1614.3845
$ws.Range("N62").Value = -5039
# Row 65: The Lumber of Their Discontent (L) / Cedar Lumber
$ws.Range("H65").Value = 2728.6843
$ws.Range("I65").Value = 2238.3845
$ws.Range("J65").Value = 3791
$ws.Range("K65").Value = 11191.9225
$ws.Range("L65").Value = 18955
$ws.Range("M65").Value = -8071.922500000001
$ws.Range("N65").Value = -25195
# Row 86: Birch, Please / Birch Lumber
$ws.Range("H86").Value = 8579
$ws.Range("I86").Value = 4159.5454
$ws.Range("J86").Value = 16681.334
$ws.Range("K86").Value = 4159.5454
$ws.Range("L86").Value = 16681.334
$ws.Range("M86").Value = -3036.5454
$ws.Range("N86").Value = -18927.334
# Row 89: Built This City on Blocks and Soul (L) / Birch Lumber
$ws.Range("H89").Value = 8579
$ws.Range("I89").Value = 4159.5454
$ws.Range("J89").Value = 16681.334
$ws.Range("K89").Value = 20797.727
$ws.Range("L89").Value = 83406.67
$ws.Range("M89").Value = -15181.727
$ws.Range("N89").Value = -94638.67

$ws = $wb.Worksheets.Item("CUL")
# Row 109: Cure for What Ails / Purple Carrot Juice
$ws.Range("H109").Value = 2181.4119
$ws.Range("I109").Value = 1438.8
$ws.Range("J109").Value = 3242.2856
$ws.Range("K109").Value = 4316.4
$ws.Range("L109").Value = 9726.856800000001
$ws.Range("M109").Value = -3276.4
$ws.Range("N109").Value = -11806.8568
# Row 113: Can't Eat Just One / Night Vinegar
$ws.Range("H113").Value = 558.0465
$ws.Range("I113").Value = 583.43475
$ws.Range("J113").Value = 528.85
$ws.Range("K113").Value = 1750.30425
$ws.Range("L113").Value = 1586.55
$ws.Range("M113").Value = 419.6957499999999
$ws.Range("N113").Value = -5926.55

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 4799.918
$ws.Range("J70").Value = 4801.087
$ws.Range("L70").Value = 4801.087
$ws.Range("N70").Value = -5341.087
# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 4799.918
$ws.Range("J73").Value = 4801.087
$ws.Range("L73").Value = 4801.087
$ws.Range("N73").Value = -6673.087

$ws = $wb.Worksheets.Item("LTW")
# Row 68: You Could Say It's a Moving Target / Wyvern Leather
$ws.Range("H68").Value = 8244.444
$ws.Range("I68").Value = 11318.182
$ws.Range("K68").Value = 11318.182
$ws.Range("M68").Value = -10569.182
# Row 71: They Call It Bloody Mary (L) / Wyvern Leather
$ws.Range("H71").Value = 8244.444
$ws.Range("I71").Value = 11318.182
$ws.Range("K71").Value = 56590.91
$ws.Range("M71").Value = -52846.91
# Row 86: Starting Off on the Wrong Foot / Serpentskin Armguards of Maiming
$ws.Range("H86").Value = 30000
$ws.Range("J86").Value = 30000
$ws.Range("L86").Value = 30000
$ws.Range("N86").Value = -32372
# Row 89: Raising the Dragoons (L) / Serpentskin Armguards of Maiming
$ws.Range("H89").Value = 30000
$ws.Range("J89").Value = 30000
$ws.Range("L89").Value = 90000
$ws.Range("N89").Value = -101856
# Row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 30780.709
$ws.Range("I93").Value = 1715.1364
$ws.Range("J93").Value = 350502
$ws.Range("K93").Value = 1715.1364
$ws.Range("L93").Value = 350502
$ws.Range("M93").Value = -467.1364000000001
$ws.Range("N93").Value = -352998
# Row 133: The Perfect Accessory / Loboskin Amulet of Fending
$ws.Range("H133").Value = 32010
$ws.Range("J133").Value = 32010
$ws.Range("L133").Value = 32010
$ws.Range("N133").Value = -37070

$ws = $wb.Worksheets.Item("WVR")
# Row 49: A Leg Up on the Cold / Linen Tights
$ws.Range("H49").Value = 6000
$ws.Range("I49").Value = 1000
$ws.Range("J49").Value = 8500
$ws.Range("K49").Value = 1000
$ws.Range("L49").Value = 8500
$ws.Range("M49").Value = -770
$ws.Range("N49").Value = -8960
# Row 62: Pride Up in Smoke / Rainbow Cloth
$ws.Range("H62").Value = 2407.3333
$ws.Range("I62").Value = 2000
$ws.Range("K62").Value = 2000
$ws.Range("M62").Value = -1376
# Row 65: Desperate for Diversionaries (L) / Rainbow Cloth
$ws.Range("H65").Value = 2407.3333
$ws.Range("I65").Value = 2000
$ws.Range("K65").Value = 10000
$ws.Range("M65").Value = -6880
